$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value, as captured from the commit diff.
$updates = [ordered]@{
    'D2' = '46.213.40'
    'E2' = '  -0.76%  '
    'D3' = '2.610.61'
    'E3' = '  +1.15%  '
    'E4' = '  +0.01%  '
    'D5' = '311.40'
    'E5' = '  +1.87%  '
    'D6' = '99.38'
    'E6' = '  -0.76%  '
    'D7' = '0.597'
    'E7' = '  -0.62%  '
    'E8' = '  +0.02%  '
    'D9' = '0.582'
    'E9' = '  +1.44%  '
    'D10' = '39.11'
    'E10' = '  +1.37%  '
    'D11' = '54.30'
    'E11' = '  -1.14%  '
    'D12' = '0.0842'
    'E12' = '  +0.30%  '
    'D13' = '8.16'
    'E13' = '  +0.04%  '
    'D14' = '2.994.61'
    'E14' = '  +0.78%  '
    'D15' = '0.107'
    'E15' = '  +1.07%  '
    'D16' = '2.591.02'
    'E16' = '  +0.46%  '
    'D17' = '0.919'
    'E17' = '  +1.40%  '
    'D18' = '14.88'
    'E18' = '  -0.17%  '
    'D19' = '46.335.87'
    'E19' = '  -0.74%  '
    'D20' = '0.0000102'
    'E20' = '  +0.78%  '
    'B21' = 'InternetComputer(DFINITY)'
    'C21' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'D21' = '12.90'
    'E21' = '  -3.85%  '
    'B22' = 'Uniswap'
    'C22' = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
    'D22' = '6.75'
    'E22' = '  +1.24%  '
    'D23' = '72.25'
    'E23' = '  +2.22%  '
    'D24' = '275.58'
    'E24' = '  +8.59%  '
    'D25' = '3.05'
    'E25' = '  +2.76%  '
    'D26' = '2.21'
    'E26' = '  +0.76%  '
    'D27' = '30.06'
    'E27' = '  +6.61%  '
    'E28' = '  -0.09%  '
    'D29' = '4.07'
    'E29' = '  +1.40%  '
    'D30' = '10.78'
    'E30' = '  +2.74%  '
    'D31' = '38.26'
    'E31' = '  -3.44%  '
    'D32' = '2.22'
    'E32' = '  -3.43%  '
    'D33' = '6.27'
    'E33' = '  +2.06%  '
    'E34' = '  -4.70%  '
    'D35' = '155.84'
    'E35' = '  +3.85%  '
    'D36' = '2.24'
    'E36' = '  -4.12%  '
    'D37' = '0.0839'
    'E37' = '  +1.16%  '
    'E38' = '  -3.97%  '
    'E39' = '  +5.37%  '
    'E40' = '  +0.74%  '
    'D41' = '23.38'
    'E41' = '  +28.24%  '
    'D42' = '15.91'
    'E42' = '  +1.16%  '
    'D43' = '0.0332'
    'E43' = '  +2.77%  '
    'D44' = '3.63'
    'E44' = '  +0.32%  '
    'D45' = '3.99'
    'E45' = '  -4.96%  '
    'D46' = '2.104.65'
    'E46' = '  +4.26%  '
    'D47' = '0.997'
    'E47' = '  -0.08%  '
    'D48' = '95.39'
    'E48' = '  +3.74%  '
    'D49' = '9.71'
    'E49' = '  +6.70%  '
    'D50' = '109.38'
    'E50' = '  +0.33%  '
    'B51' = 'Stacks'
    'C51' = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
    'D51' = '1.76'
    'E51' = '  -1.85%  '
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    # Force text storage so numeric-looking strings (e.g. "99.38") are not
    # coerced into Excel numbers -- the source data keeps these as plain text.
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
    $range.Style = "Normal"
}
